# Adds a new "2022-Q4" quarter sheet (cloned from the existing "2022-Q3"
# sheet so the formatting/styles match exactly), fills it with the new
# quarter's fund-holding data, and updates the "总计" (summary) sheet with
# the corresponding new row, shifting the rest of the historical data down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q4" sheet by cloning "2022-Q3" (so headers,
#    column A's bold/boxed style, header row style, etc. are identical),
#    then place the clone right before "2022-Q3" and rename it.
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet)

# The freshly created clone is now the active sheet, positioned directly
# before the original "2022-Q3" sheet.
$q4Sheet = $excel.ActiveSheet
$q4Sheet.Name = "2022-Q4"

# The clone has 14 data rows (rows 2-15); the new quarter only has 9, so
# drop the extra rows at the bottom before filling in the real values.
$q4Sheet.Range("A11:H15").Delete() | Out-Null

# Fund-holding rows for 2022-Q4: code, name, scale, stock position,
# position ratio, holding value (100M CNY), position rank.
$q4Data = @(
    @("014179", "中银证券远见价值混合A",           "1.59", "85.40", "4.45", "0.0708", 6),
    @("005571", "中银证券新能源灵活配置混合A",       "0.55", "89.79", "8.50", "0.0468", 1),
    @("003980", "中银证券瑞益灵活配置混合A",         "0.58", "88.28", "6.37", "0.0369", 2),
    @("005572", "中银证券新能源灵活配置混合C",       "0.26", "89.79", "8.50", "0.0221", 1),
    @("003981", "中银证券瑞益灵活配置混合C",         "0.21", "88.28", "6.37", "0.0134", 2),
    @("014180", "中银证券远见价值混合C",             "0.13", "85.40", "4.45", "0.0058", 6),
    @("011205", "兴银中证500指数增强C",              "0.66", "84.84", "0.74", "0.0049", 5),
    @("501069", "华宝标普中国Ａ股质量价值指数（LOF）", "0.14", "93.83", "2.60", "0.0036", 8),
    @("010253", "兴银中证500指数增强A",              "0.32", "84.84", "0.74", "0.0024", 5)
)

$r = 2
foreach ($row in $q4Data) {
    $q4Sheet.Cells.Item($r, 1).Value = ($r - 2)
    # Leading apostrophe forces text storage (so codes like "014179" keep
    # their leading zero and numbers like "1.59" stay text, matching the
    # other quarter sheets). The fund-name column never looks numeric, so
    # it doesn't need the apostrophe.
    $q4Sheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $q4Sheet.Cells.Item($r, 3).Value = $row[1]
    $q4Sheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $q4Sheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $q4Sheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $q4Sheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $q4Sheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: new 2022-Q4 row on top, all the
#    historical rows shift down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Grow column A's styled/boxed formatting down into the new row 10 by
# copying the last existing row before overwriting its values.
$summary.Range("A9:D9").Copy($summary.Range("A10:D10")) | Out-Null

$summaryData = @(
    @("2022-Q4", 9, 0.21),
    @("2022-Q3", 14, 0.28),
    @("2022-Q2", 68, 4.83),
    @("2022-Q1", 3, 0.12),
    @("2021-Q4", 10, 2.43),
    @("2021-Q3", 39, 5.83),
    @("2021-Q2", 27, 13.47),
    @("2021-Q1", 19, 9.14),
    @("2020-Q4", 13, 5.86)
)

$r = 2
foreach ($row in $summaryData) {
    $summary.Cells.Item($r, 1).Value = ($r - 2)
    # "20XX-QN" labels aren't numeric-looking, so they store as text as-is.
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

# Restore the originally active sheet/tab.
$summary.Activate()
